$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 221, pushing the existing rows 221-225 down to 222-226.
$ws.Rows(221).Insert()

# Populate the new row 221 with its final values (a new weekly price record).
$ws.Range("A221").Value = 10
$ws.Range("B221").Value = "Vega Modelo de Temuco"
$ws.Range("C221").Value = "La Araucanía"
$ws.Range("D221").Value = 44448
$ws.Range("E221").Value = 9
$ws.Range("F221").Value = 100112032
$ws.Range("G221").Value = "Zapallo italiano"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 270
$ws.Range("K221").Value = 17000
$ws.Range("L221").Value = 18000
$ws.Range("M221").Value = 17556
$ws.Range("N221").Value = "$/caja 60 unidades"
$ws.Range("O221").Value = "Región de Arica y Parinacota"
$ws.Range("P221").Value = 293
$ws.Range("Q221").Value = 60
$ws.Range("R221").Value = "Hortaliza"

# Match the date-time number format used by the other rows in column D.
$ws.Range("D221").NumberFormat = $ws.Range("D222").NumberFormat
